$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstract = "Background`nid=""Par1"">The outbreak of a novel coronavirus (2019-nCoV)-infected pneumonia (NCIP) is currently ongoing in China.`n`n Most of the critically ill patients received high-flow nasal cannula (HFNC) oxygen therapy.`n`n However, the experience of HFNC in this population is lacking.`n`n`nMethods`nid=""Par2"">We retrospectively screened 318 confirmed patients with NCIP in two hospitals of Chongqing, China, from January 1st to March 4th, 2020. Among them, 27 (8.4%) patients experienced severe acute respiratory failure including 17 patients (63%) treated with HFNC as first-line therapy, 9 patients (33%) treated with noninvasive ventilation (NIV) and one patient (4%) treated with invasive ventilation.`n`n HFNC failure was defined by the need of NIV or intubation as rescue therapy.`n`n`nResults`nid=""Par3"">Of the 17 HFNC patients, 7 (41%) experienced HFNC failure.`n`n The HFNC failure rate was 0% (0/6) in patients with PaO2/FiO2" + [char]0x2009 + "&gt;" + [char]0x2009 + "200" + [char]0xA0 + "mm Hg vs.`n`n 63% (7/11) in those with PaO2/FiO2" + [char]0x2009 + [char]0x2264 + [char]0x2009 + "200" + [char]0xA0 + "mm Hg (p" + [char]0x2009 + "=" + [char]0x2009 + "0.04).`n`n Compared with baseline data, the respiratory rate significantly decreased after 1" + [char]0x2013 + "2" + [char]0xA0 + "h of HFNC in successful group [median 26 (IQR: 25" + [char]0x2013 + "29) vs.`n`n 23 (22" + [char]0x2013 + "25), p" + [char]0x2009 + "=" + [char]0x2009 + "0.03].`n`n However, it did not in the unsuccessful group.`n`n After initiation of NIV as rescue therapy among the 7 patients with HFNC failure, PaO2/FiO2 significantly improved after 1" + [char]0x2013 + "2" + [char]0xA0 + "h of NIV [median 172 (150" + [char]0x2013 + "208) mmHg vs.`n`n 114 (IQR: 79" + [char]0x2013 + "130) under HFNC, p" + [char]0x2009 + "=" + [char]0x2009 + "0.04].`n`n However, two out of seven (29%) patients with NIV as rescue therapy ultimately received intubation.`n`n Among the 27 patients with severe acute respiratory failure, four patients were eventually intubated (15%).`n`n`nConclusions`nid=""Par4"">Our study indicated that HFNC was the most common ventilation support for patients with NCIP.`n`n Patients with lower PaO2/FiO2 were more likely to experience HFNC failure.`n`n`n"

$authors = "[Ke%Wang%NULL%0,  Wei%Zhao%NULL%1,  Ji%Li%NULL%1,  Weiwei%Shu%shuweiwei361@163.com%1,  Jun%Duan%duanjun412589@163.com%1]"

$ws.Range("D16").Value = $abstract
$ws.Range("E16").Value = $authors

# Setting the long, multi-line abstract text triggers the engine's implicit
# row auto-height (adds ht/customHeight to row 16). Re-run AutoFit so the
# explicit height is cleared again and the row goes back to using the
# sheet's default row height, matching the target output.
$ws.Rows.Item(16).AutoFit()
